$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13: Camila Cortes -> Camilo Cortes (also gender female -> male)
$ws.Range("A13").Value = "camilo.cortes@example.com"
$ws.Range("C13").Value = "male"
$ws.Range("F13").Value = "Camilo"

# Row 23: Javiera Sanchez -> Javier Sanchez (email + firstName only; gender stays female)
$ws.Range("A23").Value = "javier.sanchez@example.com"
$ws.Range("F23").Value = "Javier"
